$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Efna1 -> Epha5 -> sCs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha5"
$ws.Range("D2").Value = "sCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 18.57067033333334
$ws.Range("H2").Value = 55.712011
$ws.Range("I2").Value = 0.834164862818447
$ws.Range("J2").Value = 0.8341648628184472
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1979113333333334
$ws.Range("N2").Value = 0.5937340000000001
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 3.675346126563779
$ws.Range("R2").Value = 33.07811513907401
$ws.Range("S2").Value = 0.834164862818447
$ws.Range("T2").Value = 0.8341648628184472

# Row 3: FAPs -> Efna1 -> Epha5 -> sCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha5"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.333117
$ws.Range("H3").Value = 6.999351000000001
$ws.Range("I3").Value = 0.1047998907584427
$ws.Range("J3").Value = 0.1047998907584427
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1979113333333334
$ws.Range("N3").Value = 0.5937340000000001
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.4617502962926667
$ws.Range("R3").Value = 4.155752666634001
$ws.Range("S3").Value = 0.1047998907584427
$ws.Range("T3").Value = 0.1047998907584427

# Row 4: M2 -> Efna1 -> Epha5 -> sCs
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha5"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.118438
$ws.Range("H4").Value = 0.355314
$ws.Range("I4").Value = 0.005320045870673627
$ws.Range("J4").Value = 0.005320045870673628
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1979113333333334
$ws.Range("N4").Value = 0.5937340000000001
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.02344022249733334
$ws.Range("R4").Value = 0.2109620024760001
$ws.Range("S4").Value = 0.005320045870673627
$ws.Range("T4").Value = 0.005320045870673628

# Row 5: sCs -> Efna1 -> Epha5 -> sCs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha5"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.240364666666667
$ws.Range("H5").Value = 3.721094
$ws.Range("I5").Value = 0.05571520055243645
$ws.Range("J5").Value = 0.05571520055243646
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1979113333333334
$ws.Range("N5").Value = 0.5937340000000001
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.2454822249995556
$ws.Range("R5").Value = 2.209340024996
$ws.Range("S5").Value = 0.05571520055243645
$ws.Range("T5").Value = 0.05571520055243646
